# Update the "two-digit number divided by one-digit number" practice
# sheet: replace each division prompt in the tables with a new one.
# Each old prompt text is unique in the document, so a plain
# Find/Replace (MatchWholeWord, not MatchWildcards) targeting
# $d.Content is sufficient and unambiguous for each rule.
#
# NOTE: one new value ("36÷5=") happens to equal another rule's OLD
# search text, so that rule (36÷5= -> 35÷7=) is executed *before* the
# rule that produces "36÷5=" (44÷3= -> 36÷5=), to avoid the later
# write being caught by the earlier rule's search pattern.
$d = $word.ActiveDocument

$d.Content.Find.Execute("22÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷2=", 2) | Out-Null
$d.Content.Find.Execute("13÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷9=", 2) | Out-Null
$d.Content.Find.Execute("62÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷5=", 2) | Out-Null
$d.Content.Find.Execute("32÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷2=", 2) | Out-Null
$d.Content.Find.Execute("91÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷2=", 2) | Out-Null
$d.Content.Find.Execute("33÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷2=", 2) | Out-Null
$d.Content.Find.Execute("59÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷5=", 2) | Out-Null
$d.Content.Find.Execute("57÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷6=", 2) | Out-Null
$d.Content.Find.Execute("36÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷5=", 2) | Out-Null
$d.Content.Find.Execute("70÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷8=", 2) | Out-Null
$d.Content.Find.Execute("23÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷2=", 2) | Out-Null
$d.Content.Find.Execute("48÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷8=", 2) | Out-Null
$d.Content.Find.Execute("12÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷6=", 2) | Out-Null
$d.Content.Find.Execute("70÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷2=", 2) | Out-Null
$d.Content.Find.Execute("27÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷7=", 2) | Out-Null
$d.Content.Find.Execute("36÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=", 2) | Out-Null   # runs before 44÷3= rule, see note above
$d.Content.Find.Execute("44÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷5=", 2) | Out-Null
$d.Content.Find.Execute("46÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷9=", 2) | Out-Null
$d.Content.Find.Execute("29÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷9=", 2) | Out-Null
$d.Content.Find.Execute("91÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷3=", 2) | Out-Null
$d.Content.Find.Execute("47÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷3=", 2) | Out-Null
$d.Content.Find.Execute("11÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷9=", 2) | Out-Null
$d.Content.Find.Execute("10÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=", 2) | Out-Null
$d.Content.Find.Execute("78÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷3=", 2) | Out-Null
$d.Content.Find.Execute("69÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "55÷3=", 2) | Out-Null
